$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257 (pushes existing rows 257-265 down to 258-266)
$ws.Rows(257).Insert()

# Populate the newly inserted row 257 with the new weekly price record
$ws.Range("A257").Value = 6
$ws.Range("B257").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C257").Value = "Metropolitana"
$ws.Range("D257").Value = 44509
$ws.Range("E257").Value = 13
$ws.Range("F257").Value = 100112052
$ws.Range("G257").Value = "Albahaca"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 550
$ws.Range("K257").Value = 5000
$ws.Range("L257").Value = 6000
$ws.Range("M257").Value = 5545
$ws.Range("N257").Value = '$/docena de matas'
$ws.Range("O257").Value = "Región Metropolitana"
$ws.Range("P257").Value = 924
$ws.Range("Q257").Value = 6
$ws.Range("R257").Value = "Hortaliza"
